$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2457.6
$ws.Range("I11").Value = 2457.6
$ws.Range("K11").Value = 2457.6
$ws.Range("M11").Value = -2317.6
$ws.Range("H28").Value = 1146.8182
$ws.Range("I28").Value = 514.86664
$ws.Range("K28").Value = 514.86664
$ws.Range("M28").Value = -29.86663999999996
$ws.Range("H33").Value = 1309.4166
$ws.Range("I33").Value = 1471.3
$ws.Range("K33").Value = 1471.3
$ws.Range("M33").Value = -1242.3
$ws.Range("H107").Value = 20834828
$ws.Range("J107").Value = 41667260
$ws.Range("L107").Value = 41667260
$ws.Range("N107").Value = -41671100
$ws.Range("H113").Value = 64455236
$ws.Range("I113").Value = 11113777
$ws.Range("J113").Value = 100016210
$ws.Range("K113").Value = 11113777
$ws.Range("L113").Value = 100016210
$ws.Range("M113").Value = -11110523
$ws.Range("N113").Value = -100022718
$ws.Range("H133").Value = 25074996
$ws.Range("J133").Value = 25074996
$ws.Range("L133").Value = 25074996
$ws.Range("N133").Value = -25085116
$ws.Range("H138").Value = 3131254.5
$ws.Range("I138").Value = 2980.8333
$ws.Range("J138").Value = 3853164
$ws.Range("K138").Value = 8942.499899999999
$ws.Range("L138").Value = 11559492
$ws.Range("M138").Value = -3802.499899999999
$ws.Range("N138").Value = -11569772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2705481.2
$ws.Range("I32").Value = 2705481.2
$ws.Range("K32").Value = 2705481.2
$ws.Range("M32").Value = -2705194.2
$ws.Range("H132").Value = 6158.5454
$ws.Range("I132").Value = 3281.0386
$ws.Range("K132").Value = 9843.1158
$ws.Range("M132").Value = -7313.1158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 129188.875
$ws.Range("I86").Value = 203102.4
$ws.Range("J86").Value = 5999.6665
$ws.Range("K86").Value = 203102.4
$ws.Range("L86").Value = 5999.6665
$ws.Range("M86").Value = -201979.4
$ws.Range("N86").Value = -8245.666499999999
$ws.Range("H89").Value = 129188.875
$ws.Range("I89").Value = 203102.4
$ws.Range("J89").Value = 5999.6665
$ws.Range("K89").Value = 1015512
$ws.Range("L89").Value = 29998.3325
$ws.Range("M89").Value = -1009896
$ws.Range("N89").Value = -41230.3325
$ws.Range("H94").Value = 1751.1111
$ws.Range("I94").Value = 1161.0869
$ws.Range("K94").Value = 1161.0869
$ws.Range("M94").Value = -710.0869
$ws.Range("H107").Value = 4638631
$ws.Range("I107").Value = 6587529
$ws.Range("K107").Value = 6587529
$ws.Range("M107").Value = -6585609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4430.3335
$ws.Range("I16").Value = 1311.4
$ws.Range("K16").Value = 1311.4
$ws.Range("M16").Value = -1024.4
$ws.Range("H113").Value = 4430.3335
$ws.Range("I113").Value = 1311.4
$ws.Range("K113").Value = 1311.4
$ws.Range("M113").Value = 858.5999999999999
$ws.Range("H122").Value = 3884.6333
$ws.Range("I122").Value = 2587.9375
$ws.Range("K122").Value = 7763.8125
$ws.Range("M122").Value = -5313.8125
$ws.Range("H132").Value = 5902.407
$ws.Range("I132").Value = 2706.5454
$ws.Range("K132").Value = 8119.6362
$ws.Range("M132").Value = -5589.6362
$ws.Range("H134").Value = 5451.854
$ws.Range("I134").Value = 2374.476
$ws.Range("J134").Value = 7845.3706
$ws.Range("K134").Value = 7123.428
$ws.Range("L134").Value = 23536.1118
$ws.Range("M134").Value = -4588.428
$ws.Range("N134").Value = -28606.1118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5908.7095
$ws.Range("I132").Value = 2249.5
$ws.Range("J132").Value = 9811.866
$ws.Range("K132").Value = 20245.5
$ws.Range("L132").Value = 88306.79399999999
$ws.Range("M132").Value = -17715.5
$ws.Range("N132").Value = -93366.79399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1153142.9
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540
$ws.Range("H73").Value = 1153142.9
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872
$ws.Range("H122").Value = 3304798
$ws.Range("I122").Value = 9079231
$ws.Range("K122").Value = 27237693
$ws.Range("M122").Value = -27235243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5576.5835
$ws.Range("I7").Value = 3509.6365
$ws.Range("J7").Value = 7325.5386
$ws.Range("K7").Value = 3509.6365
$ws.Range("L7").Value = 7325.5386
$ws.Range("M7").Value = -3397.6365
$ws.Range("N7").Value = -7549.5386
$ws.Range("H46").Value = 4834500
$ws.Range("J46").Value = 10105727
$ws.Range("L46").Value = 10105727
$ws.Range("N46").Value = -10106103
$ws.Range("H61").Value = 3886.6287
$ws.Range("I61").Value = 2525.3914
$ws.Range("J61").Value = 6495.6665
$ws.Range("K61").Value = 2525.3914
$ws.Range("L61").Value = 6495.6665
$ws.Range("M61").Value = -2323.3914
$ws.Range("N61").Value = -6899.6665
$ws.Range("H113").Value = 3886.6287
$ws.Range("I113").Value = 2525.3914
$ws.Range("J113").Value = 6495.6665
$ws.Range("K113").Value = 2525.3914
$ws.Range("L113").Value = 6495.6665
$ws.Range("M113").Value = -355.3914
$ws.Range("N113").Value = -10835.6665
$ws.Range("H122").Value = 2952.3684
$ws.Range("I122").Value = 2080.862
$ws.Range("K122").Value = 6242.586
$ws.Range("M122").Value = -3792.586
$ws.Range("H126").Value = 5576.5835
$ws.Range("I126").Value = 3509.6365
$ws.Range("J126").Value = 7325.5386
$ws.Range("K126").Value = 10528.9095
$ws.Range("L126").Value = 21976.6158
$ws.Range("M126").Value = -8058.9095
$ws.Range("N126").Value = -26916.6158
$ws.Range("H136").Value = 14283.108
$ws.Range("I136").Value = 3533.5454
$ws.Range("K136").Value = 10600.6362
$ws.Range("M136").Value = -8050.636200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12296136
$ws.Range("I122").Value = 18002288
$ws.Range("J122").Value = 5961.4614
$ws.Range("K122").Value = 54006864
$ws.Range("L122").Value = 17884.3842
$ws.Range("M122").Value = -54004414
$ws.Range("N122").Value = -22784.3842
$ws.Range("H132").Value = 62508496
$ws.Range("I132").Value = 83342790
$ws.Range("K132").Value = 250028370
$ws.Range("M132").Value = -250025840
$ws.Range("H136").Value = 40404556
$ws.Range("I136").Value = 71429780
$ws.Range("K136").Value = 214289340
$ws.Range("M136").Value = -214286790
$ws.Range("H138").Value = 129986.336
$ws.Range("J138").Value = 129986.336
$ws.Range("L138").Value = 129986.336
$ws.Range("N138").Value = -140266.336
